$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.7047088305883165
$ws.Range("C2").Value = 0.05531215961917724
$ws.Range("D2").Value = 0.5160540542205894
$ws.Range("E2").Value = 0.166466822983935
$ws.Range("G2").Value = 0.002565318543423551
$ws.Range("I2").Value = 1.769446562148353
$ws.Range("J2").Value = 0.06393585962438664
$ws.Range("K2").Value = 0.4283934247583829
$ws.Range("L2").Value = 0.4326497763752783
$ws.Range("M2").Value = 0.2634262287300082
$ws.Range("O2").Value = 7.512710118798338
$ws.Range("B3").Value = 0.6759414607001304
$ws.Range("C3").Value = 0.05298769355981392
$ws.Range("D3").Value = 0.5147598999010512
$ws.Range("E3").Value = 0.1669860101378031
$ws.Range("G3").Value = 0.002568147267074873
$ws.Range("I3").Value = 1.779041509962653
$ws.Range("J3").Value = 0.06393023730926828
$ws.Range("K3").Value = 0.3999517337794174
$ws.Range("L3").Value = 0.4311455799168655
$ws.Range("M3").Value = 0.2581172831593364
$ws.Range("O3").Value = 7.545669195372938
$ws.Range("B4").Value = 0.658559278395245
$ws.Range("C4").Value = 0.05154171224981496
$ws.Range("D4").Value = 0.5141647700995122
$ws.Range("E4").Value = 0.1673510313855324
$ws.Range("G4").Value = 0.002569978076573484
$ws.Range("I4").Value = 1.785551347098846
$ws.Range("J4").Value = 0.06392772947632075
$ws.Range("K4").Value = 0.3826363336771266
$ws.Range("L4").Value = 0.4303803731871554
$ws.Range("M4").Value = 0.2549653554539937
$ws.Range("O4").Value = 7.568729202445525
$ws.Range("B5").Value = 0.6515471687336856
$ws.Range("C5").Value = 0.05094775534919194
$ws.Range("D5").Value = 0.5139725344370589
$ws.Range("E5").Value = 0.1675114340736137
$ws.Range("G5").Value = 0.002570747844054146
$ws.Range("I5").Value = 1.78835983004538
$ws.Range("J5").Value = 0.06392694601569548
$ws.Range("K5").Value = 0.375617817865276
$ws.Range("L5").Value = 0.4301084824219146
$ws.Range("M5").Value = 0.2537081612777641
$ws.Range("O5").Value = 7.578836861512031
$ws.Range("B6").Value = 0.6503871351341672
$ws.Range("C6").Value = 0.05084884512908161
$ws.Range("D6").Value = 0.5139436545113512
$ws.Range("E6").Value = 0.1675387734466902
$ws.Range("G6").Value = 0.002570877096760862
$ws.Range("I6").Value = 1.788835582929291
$ws.Range("J6").Value = 0.06392683035859648
$ws.Range("K6").Value = 0.3744546861184261
$ws.Range("L6").Value = 0.4300657505827061
$ws.Range("M6").Value = 0.2535010545641043
$ws.Range("O6").Value = 7.580558163427241
$ws.Range("B7").Value = 0.6584644212804562
$ws.Range("C7").Value = 0.05153372099240272
$ws.Range("D7").Value = 0.514161973778954
$ws.Range("E7").Value = 0.1673531474118519
$ws.Range("G7").Value = 0.002569988362025823
$ws.Range("I7").Value = 1.785588592784546
$ws.Range("J7").Value = 0.06392771794362684
$ws.Range("K7").Value = 0.3825415264671079
$ws.Range("L7").Value = 0.4303765445194969
$ws.Range("M7").Value = 0.2549482900082296
$ws.Range("O7").Value = 7.568862640150172
$ws.Range("B8").Value = 0.6947318400091547
$ws.Range("C8").Value = 0.05451457880499078
$ws.Range("D8").Value = 0.5155664982224692
$ws.Range("E8").Value = 0.1666362578034413
$ws.Range("G8").Value = 0.002566274427981999
$ws.Range("I8").Value = 1.772626616532158
$ws.Range("J8").Value = 0.06393372568238309
$ws.Range("K8").Value = 0.4185562941217711
$ws.Range("L8").Value = 0.4320983244575274
$ws.Range("M8").Value = 0.2615734138294705
$ws.Range("O8").Value = 7.523488880318183
$ws.Range("B9").Value = 0.7680610011598219
$ws.Range("C9").Value = 0.06021119973310363
$ws.Range("D9").Value = 0.5198988579669503
$ws.Range("E9").Value = 0.1655961962026158
$ws.Range("G9").Value = 0.002559733751965531
$ws.Range("I9").Value = 1.752109713258619
$ws.Range("J9").Value = 0.06395295652128663
$ws.Range("K9").Value = 0.4903385870037482
$ws.Range("L9").Value = 0.4367269990359972
$ws.Range("M9").Value = 0.2754153343532906
$ws.Range("O9").Value = 7.45688444170591
$ws.Range("B10").Value = 0.8232597940341293
$ws.Range("C10").Value = 0.06430605651597432
$ws.Range("D10").Value = 0.5240384370553102
$ws.Range("E10").Value = 0.1650535948620337
$ws.Range("G10").Value = 0.002555376375562239
$ws.Range("I10").Value = 1.740016392955887
$ws.Range("J10").Value = 0.06397157426810374
$ws.Range("K10").Value = 0.5437666004973494
$ws.Range("L10").Value = 0.4408860499628418
$ws.Range("M10").Value = 0.286097619734349
$ws.Range("O10").Value = 7.421558911596492
$ws.Range("B11").Value = 0.8486542066331992
$ws.Range("C11").Value = 0.06614935487024809
$ws.Range("D11").Value = 0.5261282602669866
$ws.Range("E11").Value = 0.1648545551970546
$ws.Range("G11").Value = 0.002553490429313715
$ws.Range("I11").Value = 1.735160444320556
$ws.Range("J11").Value = 0.06398100801670559
$ws.Range("K11").Value = 0.5682190860568994
$ws.Range("L11").Value = 0.4429417674448501
$ws.Range("M11").Value = 0.2910674808154923
$ws.Range("O11").Value = 7.408437304371802
$ws.Range("B12").Value = 0.8583107753023285
$ws.Range("C12").Value = 0.06684456203518607
$ws.Range("D12").Value = 0.5269492374129925
$ws.Range("E12").Value = 0.1647860310867095
$ws.Range("G12").Value = 0.002552790038278899
$ws.Range("I12").Value = 1.733414297446231
$ws.Range("J12").Value = 0.06398471806959805
$ws.Range("K12").Value = 0.5774994819125823
$ws.Range("L12").Value = 0.4437436601788107
$ws.Range("M12").Value = 0.2929651984188837
$ws.Range("O12").Value = 7.40389188522127
$ws.Range("B13").Value = 0.8562292798094404
$ws.Range("C13").Value = 0.06669496186067647
$ws.Range("D13").Value = 0.5267711104165613
$ws.Range("E13").Value = 0.1648004848089073
$ws.Range("G13").Value = 0.002552940268391652
$ws.Range("I13").Value = 1.733786240457619
$ws.Range("J13").Value = 0.0639839129347175
$ws.Range("K13").Value = 0.5754998665646838
$ws.Range("L13").Value = 0.4435699175924412
$ws.Range("M13").Value = 0.2925557932276064
$ws.Range("O13").Value = 7.404851997377136
$ws.Range("B14").Value = 0.849447855065705
$ws.Range("C14").Value = 0.06620660635343256
$ws.Range("D14").Value = 0.5261952098677796
$ws.Range("E14").Value = 0.1648487805759231
$ws.Range("G14").Value = 0.002553432531905281
$ws.Range("I14").Value = 1.735014930621311
$ws.Range("J14").Value = 0.06398131049108713
$ws.Range("K14").Value = 0.5689821763164105
$ws.Range("L14").Value = 0.443007270551476
$ws.Range("M14").Value = 0.2912232925593585
$ws.Range("O14").Value = 7.408054866270902
$ws.Range("B15").Value = 0.8452992609410046
$ws.Range("C15").Value = 0.06590710828014323
$ws.Range("D15").Value = 0.5258463062904752
$ws.Range("E15").Value = 0.1648792541967783
$ws.Range("G15").Value = 0.002553735849909786
$ws.Range("I15").Value = 1.735779607348626
$ws.Range("J15").Value = 0.0639797343220927
$ws.Range("K15").Value = 0.564992593370647
$ws.Range("L15").Value = 0.4426656819674974
$ws.Range("M15").Value = 0.2904091429246165
$ws.Range("O15").Value = 7.410071846905822
$ws.Range("B16").Value = 0.8216058759676628
$ws.Range("C16").Value = 0.06418519995804672
$ws.Range("D16").Value = 0.5239060116047654
$ws.Range("E16").Value = 0.1650675621726343
$ws.Range("G16").Value = 0.002555501557912108
$ws.Range("I16").Value = 1.740346717820081
$ws.Range("J16").Value = 0.06397097707293486
$ws.Range("K16").Value = 0.5421715059447934
$ws.Range("L16").Value = 0.44075498961314
$ws.Range("M16").Value = 0.2857750377743855
$ws.Range("O16").Value = 7.422475721748128
$ws.Range("B17").Value = 0.8071431069166977
$ws.Range("C17").Value = 0.06312387060847868
$ws.Range("D17").Value = 0.5227685625867906
$ws.Range("E17").Value = 0.1651953062659892
$ws.Range("G17").Value = 0.002556609367892356
$ws.Range("I17").Value = 1.74331370888553
$ws.Range("J17").Value = 0.06396585109328257
$ws.Range("K17").Value = 0.5282090210077399
$ws.Range("L17").Value = 0.4396247028772251
$ws.Range("M17").Value = 0.28296034885512
$ws.Range("O17").Value = 7.430839839129419
$ws.Range("B18").Value = 0.7988512918833806
$ws.Range("C18").Value = 0.06251158977354976
$ws.Range("D18").Value = 0.5221337953649225
$ws.Range("E18").Value = 0.1652732812023245
$ws.Range("G18").Value = 0.002557255613836345
$ws.Range("I18").Value = 1.745080994677231
$ws.Range("J18").Value = 0.06396299360144031
$ws.Range("K18").Value = 0.5201921108521503
$ws.Range("L18").Value = 0.4389900107399569
$ws.Range("M18").Value = 0.2813518174149294
$ws.Range("O18").Value = 7.435928181626423
$ws.Range("B19").Value = 0.7960484452111416
$ws.Range("C19").Value = 0.06230396760932422
$ws.Range("D19").Value = 0.5219222205374905
$ws.Range("E19").Value = 0.1653004557912379
$ws.Range("G19").Value = 0.002557475980098457
$ws.Range("I19").Value = 1.745689804598555
$ws.Range("J19").Value = 0.06396204173532105
$ws.Range("K19").Value = 0.5174801321660141
$ws.Range("L19").Value = 0.4387777664276626
$ws.Range("M19").Value = 0.280808987381107
$ws.Range("O19").Value = 7.437698686930958
$ws.Range("B20").Value = 0.8086799254419077
$ws.Range("C20").Value = 0.06323704063527202
$ws.Range("D20").Value = 0.5228876324464125
$ws.Range("E20").Value = 0.1651812421437242
$ws.Range("G20").Value = 0.00255649050198837
$ws.Range("I20").Value = 1.742991580895143
$ws.Range("J20").Value = 0.0639663873681755
$ws.Range("K20").Value = 0.5296939119528474
$ws.Range("L20").Value = 0.4397434288232631
$ws.Range("M20").Value = 0.2832589016317826
$ws.Range("O20").Value = 7.429920745486385
$ws.Range("B21").Value = 0.8514386362396635
$ws.Range("C21").Value = 0.06635012454096056
$ws.Range("D21").Value = 0.5263635632472017
$ws.Range("E21").Value = 0.1648344092842908
$ws.Range("G21").Value = 0.002553287568907452
$ws.Range("I21").Value = 1.734651519656715
$ws.Range("J21").Value = 0.06398207116254362
$ws.Range("K21").Value = 0.5708960200940112
$ws.Range("L21").Value = 0.4431718984172335
$ws.Range("M21").Value = 0.2916142542775617
$ws.Range("O21").Value = 7.407102617671171
$ws.Range("B22").Value = 0.8796182839772939
$ws.Range("C22").Value = 0.06836832198487741
$ws.Range("D22").Value = 0.5288077994099325
$ws.Range("E22").Value = 0.1646476373611563
$ws.Range("G22").Value = 0.002551274534514988
$ws.Range("I22").Value = 1.729741058578128
$ws.Range("J22").Value = 0.06399312358364817
$ws.Range("K22").Value = 0.5979448204508913
$ws.Range("L22").Value = 0.44554916064628
$ws.Range("M22").Value = 0.297166634220325
$ws.Range("O22").Value = 7.394657608651642
$ws.Range("B23").Value = 0.8645570337131687
$ws.Range("C23").Value = 0.06729267427945729
$ws.Range("D23").Value = 0.5274875184013297
$ws.Range("E23").Value = 0.1647436779908702
$ws.Range("G23").Value = 0.002552341605342336
$ws.Range("I23").Value = 1.732312466240273
$ws.Range("J23").Value = 0.06398715162438773
$ws.Range("K23").Value = 0.5834974696295774
$ws.Range("L23").Value = 0.4442679118535295
$ws.Range("M23").Value = 0.2941948830154999
$ws.Range("O23").Value = 7.40107408799517
$ws.Range("B24").Value = 0.8079850580254799
$ws.Range("C24").Value = 0.06318588303324191
$ws.Range("D24").Value = 0.5228337412508779
$ws.Range("E24").Value = 0.165187586411335
$ws.Range("G24").Value = 0.002556544212284149
$ws.Range("I24").Value = 1.743137023307597
$ws.Range("J24").Value = 0.06396614463955341
$ws.Range("K24").Value = 0.5290225606112244
$ws.Range("L24").Value = 0.4396897057032589
$ws.Range("M24").Value = 0.2831238957878384
$ws.Range("O24").Value = 7.43033539642218
$ws.Range("B25").Value = 0.7479894034817391
$ws.Range("C25").Value = 0.05868601748710489
$ws.Range("D25").Value = 0.5185584586656518
$ws.Range("E25").Value = 0.1658385532031446
$ws.Range("G25").Value = 0.002561424172033732
$ws.Range("I25").Value = 1.757136132662446
$ws.Range("J25").Value = 0.06394696202539496
$ws.Range("K25").Value = 0.4707973277825488
$ws.Range("L25").Value = 0.4431718984172335
$ws.Range("M25").Value = 0.2904091429246165
$ws.Range("O25").Value = 7.472510664147251
